# Applies the recorded edits to the "Artfynd" sheet:
#  - Row 11 and Row 12 swap their per-observation identity (A/Q/R/Z/AB/AC).
#  - Rows 17/18/19 rotate their per-observation identity:
#       new row17 <= old row19 data (gains Z/AB/AJ/AK/AM/AO, loses J/N/AF)
#       new row18 <= old row17 data (unchanged cell layout, only A/Q/R move)
#       new row19 <= old row18 data (loses Z/AB/AJ/AK/AM/AO, gains J/N/AF)
#
# Only the specific cells that actually change value/presence (per the
# source diff) are touched; Startdatum/Slutdatum (Y/AA) are left alone so
# their original inline text ("2023-09-26") is never round-tripped through
# a date conversion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Rows 11 <-> 12
# ---------------------------------------------------------------

$ws.Range("A11").Value = 112360935
$ws.Range("Q11").Value = 517977
$ws.Range("R11").Value = 7181358
$ws.Range("Z11").Value = "14:32"
$ws.Range("AB11").Value = "14:32"
$ws.Range("AC11").Value = "Färska och äldre ringhack"

$ws.Range("A12").Value = 112360565
$ws.Range("Q12").Value = 517939
$ws.Range("R12").Value = 7181204
$ws.Range("Z12").Value = "13:29"
$ws.Range("AB12").Value = "13:29"
$ws.Range("AC12").Value = "Både färska och äldre ringhack"

# ---------------------------------------------------------------
# Rows 17, 18, 19 rotation
# ---------------------------------------------------------------

# -- Row 17 becomes old row 19's record --
$ws.Range("A17").Value = 112382121
$ws.Range("B17").Value = 89549
$ws.Range("E17").Value = 1108
$ws.Range("F17").Value = "Harticka"
$ws.Range("G17").Value = "Pelloporus leporinus"
$ws.Range("H17").Value = "(Fr.) Krieglst."
$ws.Range("J17").ClearContents()
$ws.Range("N17").ClearContents()
$ws.Range("P17").Value = "Väster-Rissjön (Väster-Rissjön), Ås lm"
$ws.Range("Q17").Value = 517844
$ws.Range("R17").Value = 7181358
$ws.Range("Z17").Value = "13:56"
$ws.Range("AB17").Value = "13:56"
$ws.Range("AF17").ClearContents()
$ws.Range("AJ17").Value = "gran"
$ws.Range("AK17").Value = "Picea abies"
$ws.Range("AM17").Value = "Liggande död trädstam, markontakt"
$ws.Range("AO17").Value = "Horizontal, dead with ground contact # Picea abies"

# -- Row 18 becomes old row 17's record (cell layout unchanged) --
$ws.Range("A18").Value = 112375371
$ws.Range("Q18").Value = 517859
$ws.Range("R18").Value = 7181249

# -- Row 19 becomes old row 18's record --
$ws.Range("A19").Value = 112375418
$ws.Range("B19").Value = 77651
$ws.Range("E19").Value = 230405
$ws.Range("F19").Value = "Garnlav (ssp. sarmentosa)"
$ws.Range("G19").Value = "Alectoria sarmentosa subsp. sarmentosa"
$ws.Range("H19").Value = "(Ach.) Ach."
$ws.Range("J19").Value = ""
$ws.Range("N19").Value = ""
$ws.Range("P19").Value = "Väster-Rissön, Ås lm"
$ws.Range("Q19").Value = 517913
$ws.Range("R19").Value = 7181387
$ws.Range("Z19").ClearContents()
$ws.Range("AB19").ClearContents()
$ws.Range("AF19").Value = ""
$ws.Range("AJ19").ClearContents()
$ws.Range("AK19").ClearContents()
$ws.Range("AM19").ClearContents()
$ws.Range("AO19").ClearContents()
